# Fixed a bug with the trading module: append the missing trade record
# (row 5) to the BIIB noun-trade log, mirroring the formatting of the
# existing trade rows (row 4).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the previous trade row's formatting (date/number styles,
# boolean cell styles) onto the new row before writing values into it.
$ws.Range("A4:I4").Copy()
$ws.Range("A5:I5").PasteSpecial(-4122)

$ws.Range("A5").Value = 42636.593113425923   # Date
$ws.Range("B5").Value = $true                # Profitable
$ws.Range("C5").Value = 9979.36              # Principle
$ws.Range("D5").Value = 9951.5               # Start Principle
$ws.Range("E5").Value = 313.07               # BuyPrice
$ws.Range("F5").Value = 314.81               # SellPrice
$ws.Range("G5").Value = $false               # IsShortSell
$ws.Range("H5").Value = 0.56000000000000005  # Price Change %
$ws.Range("I5").Value = $false               # Strong trade
